# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G ("K") on Sheet1 is being regenerated with new strikeout (K) values
# replacing the previous Strike# counts, for rows 2-34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 5
    3  = 5
    4  = 9
    5  = 9
    6  = 4
    7  = 4
    8  = 2
    9  = 8
    10 = 8
    11 = 5
    12 = 5
    13 = 6
    14 = 7
    15 = 11
    16 = 10
    17 = 6
    18 = 6
    19 = 9
    20 = 7
    21 = 7
    22 = 10
    23 = 8
    24 = 8
    25 = 7
    26 = 9
    27 = 7
    28 = 13
    29 = 10
    30 = 7
    31 = 7
    32 = 5
    33 = 7
    34 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
